$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.603.08'
$ws.Range("E2").Value = '  -6.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.204.59'
$ws.Range("E3").Value = '  -6.94%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.96'
$ws.Range("E5").Value = '  +0.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.31'
$ws.Range("E6").Value = '  -10.95%  '

$ws.Range("E7").Value = '  -9.49%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.553'
$ws.Range("E9").Value = '  -10.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.54'
$ws.Range("E10").Value = '  -10.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.66'
$ws.Range("E11").Value = '  -3.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0822'
$ws.Range("E12").Value = '  -10.56%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.59'
$ws.Range("E13").Value = '  -11.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.105'
$ws.Range("E14").Value = '  -4.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.537.53'
$ws.Range("E15").Value = '  -7.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.848'
$ws.Range("E16").Value = '  -13.79%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.89'
$ws.Range("E17").Value = '  -10.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.190.00'
$ws.Range("E18").Value = '  -7.33%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.465.07'
$ws.Range("E19").Value = '  -6.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.91'
$ws.Range("E20").Value = '  +6.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.49'
$ws.Range("E21").Value = '  -11.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0945'
$ws.Range("E22").Value = '  -11.47%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.13'
$ws.Range("E23").Value = '  -11.46%  '

$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.17'
$ws.Range("E24").Value = '  -8.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '233.54'
$ws.Range("E25").Value = '  -10.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.10'
$ws.Range("E26").Value = '  -8.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.01'
$ws.Range("E27").Value = '  +0.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.06'
$ws.Range("E28").Value = '  -9.85%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.16'
$ws.Range("E29").Value = '  -8.51%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.43'
$ws.Range("E30").Value = '  -12.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.38'
$ws.Range("E31").Value = '  -9.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0873'
$ws.Range("E32").Value = '  -9.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '156.96'
$ws.Range("E33").Value = '  -7.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '33.73'
$ws.Range("E34").Value = '  -11.68%  '

$ws.Range("E35").Value = '  -8.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.26'
$ws.Range("E36").Value = '  +8.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.121'
$ws.Range("E37").Value = '  -7.61%  '

$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.85'
$ws.Range("E38").Value = '  +6.60%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.39'
$ws.Range("E39").Value = '  -8.77%  '

$ws.Range("E40").Value = '  -11.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.50'
$ws.Range("E41").Value = '  -11.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0317'
$ws.Range("E42").Value = '  -11.27%  '

$ws.Range("E43").Value = '  -0.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.783.78'
$ws.Range("E44").Value = '  +8.72%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '88.93'
$ws.Range("E45").Value = '  -12.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.94'
$ws.Range("E46").Value = '  -10.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.205'
$ws.Range("E47").Value = '  -12.70%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '76.68'
$ws.Range("E48").Value = '  -6.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.32'
$ws.Range("E49").Value = '  -4.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '59.91'
$ws.Range("E50").Value = '  -14.29%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.38'
$ws.Range("E51").Value = '  -10.28%  '
